# Dignity office admin removed
# Removes the "Office Admin Provider Delegate" row for "Dignity Health"
# (username "dignity.familypractice") from Sheet2, row 27.
# Deleting the entire row shifts every subsequent row up by one and lets
# Excel clean up the now-unused shared string / hyperlink relationship.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Rows.Item(27).Delete()
